$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 358, pushing existing rows 358:370 down to 359:371
$ws.Rows.Item(358).Insert()

# Populate the newly-inserted row 358 with the new record
$ws.Cells.Item(358, 1).Value = 7
$ws.Cells.Item(358, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(358, 3).Value = "Ñuble"
$ws.Cells.Item(358, 4).Value = 44509
$ws.Cells.Item(358, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(358, 5).Value = 16
$ws.Cells.Item(358, 6).Value = 100112004
$ws.Cells.Item(358, 7).Value = "Cebolla"
$ws.Cells.Item(358, 8).Value = "Sin especificar"
$ws.Cells.Item(358, 9).Value = "1a nueva(o)"
$ws.Cells.Item(358, 10).Value = 18000
$ws.Cells.Item(358, 11).Value = 950
$ws.Cells.Item(358, 12).Value = 1000
$ws.Cells.Item(358, 13).Value = 975
$ws.Cells.Item(358, 14).Value = '$/paquete 10 unidades (volumen en unidades)'
$ws.Cells.Item(358, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(358, 16).Value = 98
$ws.Cells.Item(358, 17).Value = 10
$ws.Cells.Item(358, 18).Value = "Hortaliza"
